$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark row 14 (Paint by numbers) as reserved/bought
$ws.Range("E14").Value = "Y"

# New wishlist item: LEGO Succulents 10309 (row 17)
$ws.Range("A17").Value = "LEGO Succulents 10309"
$ws.Range("A17").WrapText = $true
$ws.Range("B17").Value = "https://www.lego.com/cdn/cs/set/assets/bltd9d7b87d1b108e5e/10309_alt1.png?format=webply&fit=bounds&quality=75&width=1200&height=1200&dpr=1"
$ws.Range("C17").Value = "https://www.lego.com/en-ch/product/succulents-10309"
$ws.Range("D17").Value = "60 CHF"

# New wishlist item: LEGO Exotic Peacock 31157 (row 18)
$ws.Range("A18").Value = "LEGO Exotic Peacock 31157"
$ws.Range("A18").WrapText = $true
$ws.Range("B18").Value = "https://www.lego.com/cdn/cs/set/assets/bltfa6c9196a514bb03/31157.png?format=webply&fit=bounds&quality=75&width=1200&height=1200&dpr=1"
$ws.Range("C18").Value = "https://www.lego.com/en-ch/product/exotic-peacock-31157"
$ws.Range("D18").Value = "23 CHF"

# New wishlist item: LEGO Tiny Plants 10329 (row 19)
$ws.Range("A19").Value = "LEGO Tiny Plants 10329"
$ws.Range("A19").WrapText = $true
$ws.Range("B19").Value = "https://www.lego.com/cdn/cs/set/assets/bltb2f845ffd52a25b0/10329.png?format=webply&fit=bounds&quality=75&width=1200&height=1200&dpr=1"
$ws.Range("C19").Value = "https://www.lego.com/en-ch/product/tiny-plants-10329"
$ws.Range("D19").Value = "60 CHF"

$ws.Range("A19").Select()
